$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$newText = @"
<b><br><b>A</b>: "Let's eat Chinese food tonight." <br>
 <br>&nbsp&nbsp&nbsp&nbsp<b>B</b>: "I like Italian food more than Chinese food. "
<br><b>A</b>: "Let's go to the Samura restaurant." <br> 
 <br>&nbsp&nbsp&nbsp&nbsp<b>B</b>: "I prefer cheap restaurants." <br>
<br><b>A</b>: "Let's go to the Ying restaurant." <br> 
 <br>&nbsp&nbsp&nbsp&nbsp<b>B</b>: "Okay, I'll call to book a table." <br>
<br><b>A</b>: "Ok." <br></b>
"@

# Remove the trailing newline introduced by the here-string
$newText = $newText.TrimEnd("`r", "`n")

$ws.Range("A2").Value = $newText
